$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19: fill I19 and J19 with "X"
$ws.Range("I19").Value = "X"
$ws.Range("J19").Value = "X"

# Row 20: fill I20 and J20 with "X"
$ws.Range("I20").Value = "X"
$ws.Range("J20").Value = "X"

# Row 21: fill I21 with "M101 no existe"
$ws.Range("I21").Value = "M101 no existe"

# Update the selected cell to I21
$ws.Range("I21").Select()
